# Atomos_Profits workbook update - scheduled runner price refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) per sheet
# based on refreshed market-board data. No formulas involved - plain cached values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1928.9131
$ws.Range("I58").Value = 313.18182
$ws.Range("J58").Value = 3410
$ws.Range("K58").Value = 939.54546
$ws.Range("L58").Value = 10230
$ws.Range("M58").Value = -789.54546
$ws.Range("N58").Value = -10530
$ws.Range("H97").Value = 2519.6667
$ws.Range("J97").Value = 2519.6667
$ws.Range("L97").Value = 7559.000100000001
$ws.Range("N97").Value = -8551.000100000001
$ws.Range("H132").Value = 5885425
$ws.Range("I132").Value = 6899644
$ws.Range("K132").Value = 20698932
$ws.Range("M132").Value = -20696402

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2734.98
$ws.Range("I32").Value = 2429.1685
$ws.Range("K32").Value = 2429.1685
$ws.Range("M32").Value = -2142.1685
$ws.Range("H61").Value = 2405.9736
$ws.Range("I61").Value = 990.95
$ws.Range("J61").Value = 3978.2222
$ws.Range("K61").Value = 990.95
$ws.Range("L61").Value = 3978.2222
$ws.Range("M61").Value = -778.95
$ws.Range("N61").Value = -4402.2222
$ws.Range("H74").Value = 746.2857
$ws.Range("I74").Value = 729.4815
$ws.Range("J74").Value = 1200
$ws.Range("K74").Value = 729.4815
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = 144.5185
$ws.Range("N74").Value = -2948
$ws.Range("H77").Value = 746.2857
$ws.Range("I77").Value = 729.4815
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 3647.4075
$ws.Range("L77").Value = 6000
$ws.Range("M77").Value = 720.5925000000002
$ws.Range("N77").Value = -14736
$ws.Range("H88").Value = 2758.8572
$ws.Range("I88").Value = 2722.4
$ws.Range("J88").Value = 2850
$ws.Range("K88").Value = 2722.4
$ws.Range("L88").Value = 2850
$ws.Range("M88").Value = -2316.4
$ws.Range("N88").Value = -3662
$ws.Range("H91").Value = 2758.8572
$ws.Range("I91").Value = 2722.4
$ws.Range("J91").Value = 2850
$ws.Range("K91").Value = 2722.4
$ws.Range("L91").Value = 2850
$ws.Range("M91").Value = -1318.4
$ws.Range("N91").Value = -5658
$ws.Range("H132").Value = 22729614
$ws.Range("I132").Value = 29413324
$ws.Range("J132").Value = 5001.2
$ws.Range("K132").Value = 88239972
$ws.Range("L132").Value = 15003.6
$ws.Range("M132").Value = -88237442
$ws.Range("N132").Value = -20063.6
$ws.Range("H136").Value = 2405.9736
$ws.Range("I136").Value = 990.95
$ws.Range("J136").Value = 3978.2222
$ws.Range("K136").Value = 2972.85
$ws.Range("L136").Value = 11934.6666
$ws.Range("M136").Value = -422.8500000000004
$ws.Range("N136").Value = -17034.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2380.8708
$ws.Range("I134").Value = 1838.5769
$ws.Range("K134").Value = 5515.7307
$ws.Range("M134").Value = -2980.7307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2383630.8
$ws.Range("I31").Value = 3335624.5
$ws.Range("J31").Value = 3645.8333
$ws.Range("K31").Value = 3335624.5
$ws.Range("L31").Value = 3645.8333
$ws.Range("M31").Value = -3335329.5
$ws.Range("N31").Value = -4235.8333
$ws.Range("H34").Value = 2383630.8
$ws.Range("I34").Value = 3335624.5
$ws.Range("J34").Value = 3645.8333
$ws.Range("K34").Value = 3335624.5
$ws.Range("L34").Value = 3645.8333
$ws.Range("M34").Value = -3335422.5
$ws.Range("N34").Value = -4049.8333
$ws.Range("H38").Value = 7400
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H46").Value = 7400
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H58").Value = 19234530
$ws.Range("I58").Value = 2303.6667
$ws.Range("J58").Value = 45460292
$ws.Range("K58").Value = 2303.6667
$ws.Range("L58").Value = 45460292
$ws.Range("M58").Value = -2100.6667
$ws.Range("N58").Value = -45460698
$ws.Range("H62").Value = 3499
$ws.Range("I62").Value = 2760.5
$ws.Range("J62").Value = 4237.5
$ws.Range("K62").Value = 2760.5
$ws.Range("L62").Value = 4237.5
$ws.Range("M62").Value = -2136.5
$ws.Range("N62").Value = -5485.5
$ws.Range("H65").Value = 3499
$ws.Range("I65").Value = 2760.5
$ws.Range("J65").Value = 4237.5
$ws.Range("K65").Value = 13802.5
$ws.Range("L65").Value = 21187.5
$ws.Range("M65").Value = -10682.5
$ws.Range("N65").Value = -27427.5
$ws.Range("H122").Value = 3046.611
$ws.Range("J122").Value = 3402.1428
$ws.Range("L122").Value = 10206.4284
$ws.Range("N122").Value = -15106.4284
$ws.Range("H136").Value = 19234530
$ws.Range("I136").Value = 2303.6667
$ws.Range("J136").Value = 45460292
$ws.Range("K136").Value = 6911.000100000001
$ws.Range("L136").Value = 136380876
$ws.Range("M136").Value = -4361.000100000001
$ws.Range("N136").Value = -136385976

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 117.32
$ws.Range("I12").Value = 34
$ws.Range("J12").Value = 133.19048
$ws.Range("K12").Value = 102
$ws.Range("L12").Value = 399.5714400000001
$ws.Range("M12").Value = 71
$ws.Range("N12").Value = -745.5714400000001
$ws.Range("H87").Value = 11458.667
$ws.Range("I87").Value = 5996.6665
$ws.Range("J87").Value = 15100
$ws.Range("K87").Value = 17989.9995
$ws.Range("L87").Value = 45300
$ws.Range("M87").Value = -16741.9995
$ws.Range("N87").Value = -47796
$ws.Range("H90").Value = 11458.667
$ws.Range("I90").Value = 5996.6665
$ws.Range("J90").Value = 15100
$ws.Range("K90").Value = 53969.9985
$ws.Range("L90").Value = 135900
$ws.Range("M90").Value = -47729.9985
$ws.Range("N90").Value = -148380
$ws.Range("H101").Value = 3972.7144
$ws.Range("J101").Value = 3972.7144
$ws.Range("L101").Value = 11918.1432
$ws.Range("N101").Value = -16786.1432
$ws.Range("H107").Value = 1372.1578
$ws.Range("I107").Value = 912.5
$ws.Range("J107").Value = 1706.4546
$ws.Range("K107").Value = 2737.5
$ws.Range("L107").Value = 5119.3638
$ws.Range("M107").Value = -817.5
$ws.Range("N107").Value = -8959.363799999999
$ws.Range("H120").Value = 17869.111
$ws.Range("H123").Value = 1199.7059
$ws.Range("I123").Value = 932.8570999999999
$ws.Range("J123").Value = 1630.7693
$ws.Range("K123").Value = 2798.5713
$ws.Range("L123").Value = 4892.3079
$ws.Range("M123").Value = -348.5712999999996
$ws.Range("N123").Value = -9792.3079
$ws.Range("H124").Value = 15951.429
$ws.Range("J124").Value = 50950
$ws.Range("L124").Value = 152850
$ws.Range("N124").Value = -162670
$ws.Range("H125").Value = 1970
$ws.Range("J125").Value = 2996.6667
$ws.Range("L125").Value = 8990.000100000001
$ws.Range("N125").Value = -18830.0001
$ws.Range("H126").Value = 1776.25
$ws.Range("I126").Value = 1076.6666
$ws.Range("J126").Value = 2196
$ws.Range("K126").Value = 3229.9998
$ws.Range("L126").Value = 6588
$ws.Range("M126").Value = 1710.0002
$ws.Range("N126").Value = -16468
$ws.Range("H130").Value = 1899.7778
$ws.Range("J130").Value = 2228.2856
$ws.Range("L130").Value = 6684.8568
$ws.Range("N130").Value = -16724.8568
$ws.Range("H136").Value = 2347.5264
$ws.Range("I136").Value = 1606.1538
$ws.Range("J136").Value = 3953.8333
$ws.Range("K136").Value = 4818.4614
$ws.Range("L136").Value = 11861.4999
$ws.Range("M136").Value = 281.5385999999999
$ws.Range("N136").Value = -22061.4999
$ws.Range("H138").Value = 1434.2667
$ws.Range("I138").Value = 849.1
$ws.Range("K138").Value = 2547.3
$ws.Range("M138").Value = 2592.7
$ws.Range("H139").Value = 9293.611000000001
$ws.Range("I139").Value = 2890
$ws.Range("J139").Value = 14416.5
$ws.Range("K139").Value = 8670
$ws.Range("L139").Value = 43249.5
$ws.Range("M139").Value = -3530
$ws.Range("N139").Value = -53529.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4938
$ws.Range("I70").Value = 5057.143
$ws.Range("J70").Value = 4660
$ws.Range("K70").Value = 5057.143
$ws.Range("L70").Value = 4660
$ws.Range("M70").Value = -4787.143
$ws.Range("N70").Value = -5200
$ws.Range("H73").Value = 4938
$ws.Range("I73").Value = 5057.143
$ws.Range("J73").Value = 4660
$ws.Range("K73").Value = 5057.143
$ws.Range("L73").Value = 4660
$ws.Range("M73").Value = -4121.143
$ws.Range("N73").Value = -6532
$ws.Range("H111").Value = 17733.334
$ws.Range("J111").Value = 17733.334
$ws.Range("L111").Value = 17733.334
$ws.Range("N111").Value = -23867.334
$ws.Range("H132").Value = 3422.3333
$ws.Range("J132").Value = 3739.5
$ws.Range("L132").Value = 11218.5
$ws.Range("N132").Value = -16278.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3277.5715
$ws.Range("I132").Value = 1862.909
$ws.Range("J132").Value = 4192.9414
$ws.Range("K132").Value = 5588.727000000001
$ws.Range("L132").Value = 12578.8242
$ws.Range("M132").Value = -3058.727000000001
$ws.Range("N132").Value = -17638.8242

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2850
$ws.Range("I107").Value = 200
$ws.Range("J107").Value = 3733.3333
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 11199.9999
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -15039.9999
$ws.Range("H122").Value = 717339.5
$ws.Range("J122").Value = 5333
$ws.Range("L122").Value = 15999
$ws.Range("N122").Value = -20899
$ws.Range("H132").Value = 192535.75
$ws.Range("I132").Value = 271747.34
$ws.Range("K132").Value = 815242.02
$ws.Range("M132").Value = -812712.02
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
